$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.807.34"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.644.84"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").Value = "'216.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'19.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("D11").Value = "'0.0840"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "1.869.37"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "1.647.49"
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "'64.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("D17").Value = "26.794.83"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").Value = "'214.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "'4.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'2.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.33%  "
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D25").Value = "'144.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "'7.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("D34").Value = "1.291.92"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("E37").Value = "  -5.72%  "
$ws.Range("D38").Value = "'0.539"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").Value = "1.795.95"
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("D45").Value = "'60.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("D46").Value = "'91.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0520"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0978"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
